$d = $word.ActiveDocument

# Step 1: remove paragraphs 1-4 (old title block)
$p1 = $d.Paragraphs(1)
$p4 = $d.Paragraphs(4)
$rngDel = $d.Range($p1.Range.Start, $p4.Range.End)
$rngDel.Delete()

# Step 2: remove the _GoBack bookmark
try { $d.Bookmarks("_GoBack").Delete() } catch {}

# Step 3: insert new table + heading paragraphs before the remaining paragraph
$target = $d.Paragraphs(1).Range
$target.Collapse(1)
$xml = '<w:tbl>
      <w:tblPr>
        <w:tblW w:w="9356" w:type="dxa"/>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="4053"/>
        <w:gridCol w:w="5303"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4053" w:type="dxa"/>
            <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:firstLine="0"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>BAN CÔNG TÁC XDĐSVH - ĐTVM</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:firstLine="0"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>PHƯỜNG/XÃ</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>………………</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:firstLine="0"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>BVĐ XDĐSVH – ĐTVM KHÓM</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>/ẤP</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>…..</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:firstLine="0"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:i/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:firstLine="0"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:i/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="5303" w:type="dxa"/>
            <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:firstLine="0"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>CỘNG HÒA XÃ HỘI CHỦ NGHĨA VIỆT NAM</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:firstLine="0"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Độc lập – Tự do – Hạnh phúc</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:firstLine="0"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <mc:AlternateContent>
                <mc:Choice Requires="wps">
                  <w:drawing>
                    <wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251659264" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="3AB30684" wp14:editId="010DCA7A">
                      <wp:simplePos x="0" y="0"/>
                      <wp:positionH relativeFrom="column">
                        <wp:posOffset>779145</wp:posOffset>
                      </wp:positionH>
                      <wp:positionV relativeFrom="paragraph">
                        <wp:posOffset>45085</wp:posOffset>
                      </wp:positionV>
                      <wp:extent cx="1934210" cy="0"/>
                      <wp:effectExtent l="13970" t="5715" r="13970" b="13335"/>
                      <wp:wrapNone/>
                      <wp:docPr id="1852301150" name="Straight Arrow Connector 2"/>
                      <wp:cNvGraphicFramePr>
                        <a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"/>
                      </wp:cNvGraphicFramePr>
                      <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
                        <a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape">
                          <wps:wsp>
                            <wps:cNvCnPr>
                              <a:cxnSpLocks noChangeShapeType="1"/>
                            </wps:cNvCnPr>
                            <wps:spPr bwMode="auto">
                              <a:xfrm>
                                <a:off x="0" y="0"/>
                                <a:ext cx="1934210" cy="0"/>
                              </a:xfrm>
                              <a:prstGeom prst="straightConnector1">
                                <a:avLst/>
                              </a:prstGeom>
                              <a:noFill/>
                              <a:ln w="9525">
                                <a:solidFill>
                                  <a:srgbClr val="000000"/>
                                </a:solidFill>
                                <a:round/>
                                <a:headEnd/>
                                <a:tailEnd/>
                              </a:ln>
                              <a:extLst>
                                <a:ext uri="{909E8E84-426E-40DD-AFC4-6F175D3DCCD1}">
                                  <a14:hiddenFill xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main">
                                    <a:noFill/>
                                  </a14:hiddenFill>
                                </a:ext>
                              </a:extLst>
                            </wps:spPr>
                            <wps:bodyPr/>
                          </wps:wsp>
                        </a:graphicData>
                      </a:graphic>
                      <wp14:sizeRelH relativeFrom="page">
                        <wp14:pctWidth>0</wp14:pctWidth>
                      </wp14:sizeRelH>
                      <wp14:sizeRelV relativeFrom="page">
                        <wp14:pctHeight>0</wp14:pctHeight>
                      </wp14:sizeRelV>
                    </wp:anchor>
                  </w:drawing>
                </mc:Choice>
                <mc:Fallback>
                  <w:pict>
                    <v:shapetype w14:anchorId="40A68B45" id="_x0000_t32" coordsize="21600,21600" o:spt="32" o:oned="t" path="m,l21600,21600e" filled="f">
                      <v:path arrowok="t" fillok="f" o:connecttype="none"/>
                      <o:lock v:ext="edit" shapetype="t"/>
                    </v:shapetype>
                    <v:shape id="Straight Arrow Connector 2" o:spid="_x0000_s1026" type="#_x0000_t32" style="position:absolute;margin-left:61.35pt;margin-top:3.55pt;width:152.3pt;height:0;z-index:251659264;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-height-percent:0;mso-width-relative:page;mso-height-relative:page" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQD26toAuAEAAFYDAAAOAAAAZHJzL2Uyb0RvYy54bWysU8Fu2zAMvQ/YPwi6L46zdViNOD2k6y7d&#10;FqDdBzCSbAuTRYFUYufvJ6lJVmy3YT4IlEg+Pj7S67t5dOJoiC36VtaLpRTGK9TW96388fzw7pMU&#10;HMFrcOhNK0+G5d3m7Zv1FBqzwgGdNiQSiOdmCq0cYgxNVbEazAi8wGB8cnZII8R0pb7SBFNCH121&#10;Wi4/VhOSDoTKMKfX+xen3BT8rjMqfu86NlG4ViZusZxUzn0+q80amp4gDFadacA/sBjB+lT0CnUP&#10;EcSB7F9Qo1WEjF1cKBwr7DqrTOkhdVMv/+jmaYBgSi9JHA5Xmfj/wapvx63fUaauZv8UHlH9ZOFx&#10;O4DvTSHwfAppcHWWqpoCN9eUfOGwI7GfvqJOMXCIWFSYOxozZOpPzEXs01VsM0eh0mN9+/7Dqk4z&#10;URdfBc0lMRDHLwZHkY1WciSw/RC36H0aKVJdysDxkWOmBc0lIVf1+GCdK5N1XkytvL1Z3ZQERmd1&#10;duYwpn6/dSSOkHejfKXH5HkdRnjwuoANBvTnsx3Buhc7FXf+LE1WI68eN3vUpx1dJEvDKyzPi5a3&#10;4/W9ZP/+HTa/AAAA//8DAFBLAwQUAAYACAAAACEAJ/5d/NsAAAAHAQAADwAAAGRycy9kb3ducmV2&#10;LnhtbEyOTU/DMBBE70j8B2uRuKDWiflICXGqCokDR9pKXN14SQLxOoqdJvTXs3Apx6cZzbxiPbtO&#10;HHEIrScN6TIBgVR521KtYb97WaxAhGjIms4TavjGAOvy8qIwufUTveFxG2vBIxRyo6GJsc+lDFWD&#10;zoSl75E4+/CDM5FxqKUdzMTjrpMqSR6kMy3xQ2N6fG6w+tqOTgOG8T5NNo+u3r+eppt3dfqc+p3W&#10;11fz5glExDmey/Crz+pQstPBj2SD6JiVyriqIUtBcH6nslsQhz+WZSH/+5c/AAAA//8DAFBLAQIt&#10;ABQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAAAAAAAABbQ29udGVudF9UeXBlc10u&#10;eG1sUEsBAi0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAAAAAAAAAAAAAALwEAAF9yZWxzLy5y&#10;ZWxzUEsBAi0AFAAGAAgAAAAhAPbq2gC4AQAAVgMAAA4AAAAAAAAAAAAAAAAALgIAAGRycy9lMm9E&#10;b2MueG1sUEsBAi0AFAAGAAgAAAAhACf+XfzbAAAABwEAAA8AAAAAAAAAAAAAAAAAEgQAAGRycy9k&#10;b3ducmV2LnhtbFBLBQYAAAAABAAEAPMAAAAaBQAAAAA=&#10;"/>
                  </w:pict>
                </mc:Fallback>
              </mc:AlternateContent>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:firstLine="0"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:firstLine="0"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:i/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:i/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>……………, ngày    tháng   năm 20…..</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:firstLine="0"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
                <w:b/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:ind w:firstLine="0"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:b/>
          <w:sz w:val="6"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:b/>
          <w:noProof/>
          <w:szCs w:val="28"/>
          <w:lang w:val="vi-VN" w:eastAsia="vi-VN"/>
        </w:rPr>
        <mc:AlternateContent>
          <mc:Choice Requires="wps">
            <w:drawing>
              <wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251660288" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="091D599B" wp14:editId="46033406">
                <wp:simplePos x="0" y="0"/>
                <wp:positionH relativeFrom="column">
                  <wp:posOffset>5505450</wp:posOffset>
                </wp:positionH>
                <wp:positionV relativeFrom="paragraph">
                  <wp:posOffset>-1680209</wp:posOffset>
                </wp:positionV>
                <wp:extent cx="752475" cy="434340"/>
                <wp:effectExtent l="0" t="0" r="28575" b="22860"/>
                <wp:wrapNone/>
                <wp:docPr id="996151731" name="Rectangle 1"/>
                <wp:cNvGraphicFramePr>
                  <a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"/>
                </wp:cNvGraphicFramePr>
                <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
                  <a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape">
                    <wps:wsp>
                      <wps:cNvSpPr>
                        <a:spLocks noChangeArrowheads="1"/>
                      </wps:cNvSpPr>
                      <wps:spPr bwMode="auto">
                        <a:xfrm>
                          <a:off x="0" y="0"/>
                          <a:ext cx="752475" cy="434340"/>
                        </a:xfrm>
                        <a:prstGeom prst="rect">
                          <a:avLst/>
                        </a:prstGeom>
                        <a:solidFill>
                          <a:srgbClr val="FFFFFF"/>
                        </a:solidFill>
                        <a:ln w="9525">
                          <a:solidFill>
                            <a:srgbClr val="000000"/>
                          </a:solidFill>
                          <a:miter lim="800000"/>
                          <a:headEnd/>
                          <a:tailEnd/>
                        </a:ln>
                      </wps:spPr>
                      <wps:txbx>
                        <w:txbxContent>
                          <w:p>
                            <w:pPr>
                              <w:ind w:firstLine="0"/>
                              <w:jc w:val="center"/>
                            </w:pPr>
                            <w:r>
                              <w:t>Mẫu 02</w:t>
                            </w:r>
                          </w:p>
                        </w:txbxContent>
                      </wps:txbx>
                      <wps:bodyPr rot="0" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" anchor="t" anchorCtr="0" upright="1">
                        <a:noAutofit/>
                      </wps:bodyPr>
                    </wps:wsp>
                  </a:graphicData>
                </a:graphic>
                <wp14:sizeRelH relativeFrom="page">
                  <wp14:pctWidth>0</wp14:pctWidth>
                </wp14:sizeRelH>
                <wp14:sizeRelV relativeFrom="page">
                  <wp14:pctHeight>0</wp14:pctHeight>
                </wp14:sizeRelV>
              </wp:anchor>
            </w:drawing>
          </mc:Choice>
          <mc:Fallback>
            <w:pict>
              <v:rect w14:anchorId="091D599B" id="Rectangle 1" o:spid="_x0000_s1026" style="position:absolute;left:0;text-align:left;margin-left:433.5pt;margin-top:-132.3pt;width:59.25pt;height:34.2pt;z-index:251660288;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-height-percent:0;mso-width-relative:page;mso-height-relative:page;v-text-anchor:top" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQABPFmsEAIAACAEAAAOAAAAZHJzL2Uyb0RvYy54bWysU9tu2zAMfR+wfxD0vjjJkrU14hRFugwD&#10;ugvQ7QNkWbaFSaJGKbGzrx+lpGl2eRomAwJpUkeHh9TqdrSG7RUGDa7is8mUM+UkNNp1Ff/6Zfvq&#10;mrMQhWuEAacqflCB365fvlgNvlRz6ME0ChmBuFAOvuJ9jL4siiB7ZUWYgFeOgi2gFZFc7IoGxUDo&#10;1hTz6fRNMQA2HkGqEOjv/THI1xm/bZWMn9o2qMhMxYlbzDvmvU57sV6JskPhey1PNMQ/sLBCO7r0&#10;DHUvomA71H9AWS0RArRxIsEW0LZaqlwDVTOb/lbNYy+8yrWQOMGfZQr/D1Z+3D/6z5ioB/8A8ltg&#10;Dja9cJ26Q4ShV6Kh62ZJqGLwoTwfSE6go6wePkBDrRW7CFmDsUWbAKk6NmapD2ep1RiZpJ9Xy/ni&#10;asmZpNDiNX25FYUonw57DPGdAsuSUXGkTmZwsX8IMZER5VNKJg9GN1ttTHawqzcG2V5Q17d5Zf5U&#10;42WacWyo+M1yvszIv8TCJcQ0r79BWB1pfI22Fb8+J4kyqfbWNXm4otDmaBNl404yJuXSkIYyjvVI&#10;icmsoTmQoAjHMaVnRUYP+IOzgUa04uH7TqDizLx31JSb2YJkYzE7i+XVnBy8jNSXEeEkQVU8cnY0&#10;N/H4DnYeddfTTbMsg4M7amSrs8jPrE68aQyz9qcnk+b80s9Zzw97/RMAAP//AwBQSwMEFAAGAAgA&#10;AAAhABZxybXjAAAADQEAAA8AAABkcnMvZG93bnJldi54bWxMj8FOwzAQRO9I/IO1SNxap4GaJI1T&#10;IVCROLbphZsTmyQlXkex0wa+nuVUjrMzmn2Tb2fbs7MZfedQwmoZATNYO91hI+FY7hYJMB8UatU7&#10;NBK+jYdtcXuTq0y7C+7N+RAaRiXoMyWhDWHIOPd1a6zySzcYJO/TjVYFkmPD9aguVG57HkeR4FZ1&#10;SB9aNZiX1tRfh8lKqLr4qH725Vtk091DeJ/L0/TxKuX93fy8ARbMHK5h+MMndCiIqXITas96CYl4&#10;oi1BwiIWjwIYRdJkvQZW0WmVihh4kfP/K4pfAAAA//8DAFBLAQItABQABgAIAAAAIQC2gziS/gAA&#10;AOEBAAATAAAAAAAAAAAAAAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAh&#10;ADj9If/WAAAAlAEAAAsAAAAAAAAAAAAAAAAALwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAh&#10;AAE8WawQAgAAIAQAAA4AAAAAAAAAAAAAAAAALgIAAGRycy9lMm9Eb2MueG1sUEsBAi0AFAAGAAgA&#10;AAAhABZxybXjAAAADQEAAA8AAAAAAAAAAAAAAAAAagQAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAA&#10;BAAEAPMAAAB6BQAAAAA=&#10;">
                <v:textbox>
                  <w:txbxContent>
                    <w:p>
                      <w:pPr>
                        <w:ind w:firstLine="0"/>
                        <w:jc w:val="center"/>
                      </w:pPr>
                      <w:r>
                        <w:t>Mẫu 02</w:t>
                      </w:r>
                    </w:p>
                  </w:txbxContent>
                </v:textbox>
              </v:rect>
            </w:pict>
          </mc:Fallback>
        </mc:AlternateContent>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="720"/>
          <w:tab w:val="left" w:pos="2565"/>
        </w:tabs>
        <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
        <w:ind w:firstLine="0"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:b/>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:b/>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>BẢNG TỰ ĐÁNH GIÁ</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="720"/>
          <w:tab w:val="left" w:pos="2565"/>
        </w:tabs>
        <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
        <w:ind w:firstLine="0"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:b/>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>Thực hiện các Tiêu chuẩn xây dựng “Khóm văn minh đô thị” năm ……..</w:t>
      </w:r>
    </w:p>
    '
$target.InsertXML($xml) | Out-Null

Write-Host "Paragraphs:" $d.Paragraphs.Count
Write-Host "Tables:" $d.Tables.Count
